# ADD results from server
# Updates a handful of numeric result cells (row 2) on sheets "2025",
# "2030" and "2035" with refreshed values received from the server.

$wb = $excel.ActiveWorkbook

# ----- Sheet "2025" -----
$ws = $wb.Worksheets.Item("2025")
$ws.Range("E2").Value = 0.380801491414829
$ws.Range("G2").Value = 0.2494892361374987
$ws.Range("I2").Value = 0.3490586900963446
$ws.Range("L2").Value = 0.620926
$ws.Range("M2").Value = 0.07877716666666652
$ws.Range("N2").Value = 12.58268103604887
$ws.Range("O2").Value = 3.039902505957998

# ----- Sheet "2030" -----
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 0.06049637743582365
$ws.Range("E2").Value = 0.3791936161318825
$ws.Range("I2").Value = 0.3842421432369887
$ws.Range("L2").Value = 0.4679873757707057
$ws.Range("M2").Value = 0.07840229161376959
$ws.Range("N2").Value = 9.430249752715472
$ws.Range("O2").Value = 4.051469887508581

# ----- Sheet "2035" -----
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.06937835746130364
$ws.Range("B2").Value = 0.02700184948922742
$ws.Range("E2").Value = 0.1902207139172843
$ws.Range("I2").Value = 0.4643640683185117
$ws.Range("M2").Value = 0.03858787505289721
$ws.Range("N2").Value = 9.002099109637758
$ws.Range("O2").Value = 5.364149887256674
